# Cosenza.xlsx update: extend each of the 5 sheets (Nuovi casi, Deceduti,
# Dimessi Guariti, Ricoveri, Terapia) with new daily data.
#
# - Rows 476-484: fill in column C (daily count) and column D (trailing
#   7-day AVERAGE formula, continuing the existing moving-average series).
# - Rows 476-510: column A already/now holds the date serials 44374..44408
#   (2021-07-15 .. 2021-08-18 in the sheet's own numbering scheme); rows
#   485-510 only get a date, no C/D data yet (mirrors upstream source).
# - Select A476:D484 on every sheet (matches the new selection left by the
#   author after pasting the data in), ending on "Terapia" so it stays the
#   active tab exactly as before the edit.

$wb = $excel.ActiveWorkbook

$wsNuoviCasi = $wb.Worksheets.Item("Nuovi casi")
$wsDeceduti = $wb.Worksheets.Item("Deceduti")
$wsGuariti = $wb.Worksheets.Item("Dimessi   Guariti")
$wsRicoveri = $wb.Worksheets.Item("Ricoveri")
$wsTerapia = $wb.Worksheets.Item("Terapia")

# New daily values for column C, rows 476..484 (one array per sheet).
$cNuoviCasi = @(17, 1, 7, 11, 24, 4, 41, 7, 9)
$cDeceduti = @(0, 0, 0, 0, 0, 0, 0, 1, 0)
$cGuariti = @(156, 8, 130, 140, 105, 264, 1304, 225, 242)
$cRicoveri = @(24, 25, 25, 25, 20, 17, 19, 18, 24)
$cTerapia = @(3, 3, 3, 3, 4, 4, 4, 4, 4)

function Fill-Sheet {
    param($ws, $cvals)

    # Column C + D for rows 476..484 (continuing the existing
    # AVERAGE(C[n-6]:C[n]) trailing 7-day moving average).
    for ($i = 0; $i -lt $cvals.Length; $i++) {
        $r = 476 + $i
        $ws.Cells.Item($r, 3).Value = $cvals[$i]
        $ws.Cells.Item($r, 4).Formula = "=AVERAGE(C" + ($r - 6) + ":C" + $r + ")"
    }

    # Match the number formatting of the preceding rows (copy formats only).
    $ws.Range("C475").Copy()
    $ws.Range("C476:C484").PasteSpecial(-4122)
    $ws.Range("D475").Copy()
    $ws.Range("D476:D484").PasteSpecial(-4122)

    # Extend column A with the date series through row 510 (476-479 already
    # carry their dates; 480-510 are brand new rows).
    for ($r = 480; $r -le 510; $r++) {
        $ws.Cells.Item($r, 1).Value = 44374 + ($r - 476)
    }

    $ws.Activate()
    $ws.Range("A476:D484").Select()
}

Fill-Sheet $wsNuoviCasi $cNuoviCasi
Fill-Sheet $wsDeceduti $cDeceduti
Fill-Sheet $wsGuariti $cGuariti
Fill-Sheet $wsRicoveri $cRicoveri
Fill-Sheet $wsTerapia $cTerapia
